# Add a new RA ("XXXXX") as a new list item after "Luigi Caloi",
# reusing the same ListParagraph/numPr formatting, and move the
# "_GoBack" bookmark so it sits right after the newly typed text.

$d = $word.ActiveDocument
$newName = "XXXXX"

# The last paragraph in the document is the "Luigi Caloi" list item.
# Inserting a paragraph after it clones its paragraph formatting
# (style + numbering), which is what we want for the new entry.
$last = $d.Paragraphs.Last
$last.Range.InsertParagraphAfter()

# Type the new RA's name followed by a one-character sentinel so the
# bookmark can be anchored immediately after the real text without the
# (buggy, for this engine) case of landing exactly on the paragraph mark.
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = $newName + "Z"

$newPara = $d.Paragraphs.Last
$start = $newPara.Range.Start
$afterName = $start + $newName.Length
$bmRange = $d.Range($afterName, $afterName)

# Re-adding "_GoBack" automatically removes it from its previous
# location (the end of the "Luigi Caloi" paragraph) and places it here.
$d.Bookmarks.Add("_GoBack", $bmRange)

# Remove the sentinel character now that the bookmark is anchored.
$sentinel = $d.Range($afterName, $afterName + 1)
$sentinel.Delete() | Out-Null
